$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update correct/total marks on the marksheet (row 11 "Marking", row 12 "Total")
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 80
$ws.Range("E12").Value = "80/140"
